$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: search forward for literal text starting at/after $cursor.End and
# return a Range covering the match. Throws if not found (fail fast/loud).
# ---------------------------------------------------------------------------
function Step-Find($cursor, [string]$text) {
    $r = $d.Range($cursor.End, $d.Content.End)
    $ok = $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Step-Find: text not found: $text"
    }
    return $r
}

# Replace the text of a found range, then return a cursor collapsed to the
# end of the newly-inserted text (so the next search starts right after it).
function Replace-Text($cursor, [string]$oldText, [string]$newText) {
    $r = Step-Find $cursor $oldText
    $r.Text = $newText
    $r.Collapse(0)
    return $r
}

# Delete the found range entirely (used where whole runs disappear), and
# return a cursor collapsed at the deletion point.
function Delete-Text($cursor, [string]$oldText) {
    $r = Step-Find $cursor $oldText
    $r.Delete()
    $r.Collapse(0)
    return $r
}

# =============================================================================
# Title / byline / email
# =============================================================================
$cur = $d.Range(0, 0)
$cur = Replace-Text $cur "Quantum Technology in Biology: Unveiling Potential" "Understanding Our Governance: An Introduction to High School Civics"
$cur = Replace-Text $cur "Olivia Greenleaf" "Alex Thompson"
$cur = Replace-Text $cur "olivia" "alex"
$cur = Replace-Text $cur "greenleaf@quantum-bio" "thompson@hsed"
$cur = Replace-Text $cur "org" "edu"

# =============================================================================
# Body paragraph 1 (two manual line breaks inside it)
# =============================================================================
$cur = Replace-Text $cur "Stand at the intersection of biology and quantum technology, where nature's secrets meld with the enigmatic world of the quantum" "Delving into the intricate world of governance is an essential endeavor for high school students as they prepare for their roles as active and informed citizens"

$cur = Replace-Text $cur " This novel frontier offers unprecedented insights into the intricate mechanisms that govern life, promising innovations that could revolutionize medicine and biology" " This essay aims to provide a comprehensive overview of civics, exploring fundamental concepts, structures, and processes that shape our political systems"

$cur = Delete-Text $cur " Quantum phenomena such as superposition and entanglement hold immense potential for unraveling the mysteries of cells, molecules, and genetic processes."

$cur = Replace-Text $cur " With meticulous experimentation and theoretical modeling, scientists embark on a quest to decipher the language of life at its most fundamental level" " We will embark on a journey through the foundational principles of democracy, the complexities of government branches, and the dynamic relationship between citizens and their leaders"

$cur = Replace-Text $cur "From intricate nanoscale devices that probe biological systems with extreme precision to advanced imaging techniques that reveal the dynamic choreography of molecules, quantum technology is redefining our understanding of life" "Governments, in their myriad forms, serve as the backbone of organized societies, establishing rules, regulations, and institutions to maintain order, protect rights, and facilitate collective decision-making"

$cur = Replace-Text $cur " Researchers explore the possibility of harnessing quantum effects to design targeted therapies, harness cellular processes for novel materials, and unravel the secrets of DNA replication" " As we navigate this intricate landscape, we will unravel the delicate balance between individual liberties and collective responsibilities, examining how citizens can actively engage in shaping their governance"

$cur = Delete-Text $cur " The journey into quantum biology is fraught with challenges, requiring a delicate balance between quantum coherence and the complexities of living organisms."

$cur = Replace-Text $cur " Yet, with each breakthrough experiment, we approach a deeper comprehension of life's enigmatic dance" " Additionally, we will delve into the historical evolution of governance systems, tracing the transformative shifts from ancient democracies to modern representative republics"

$cur = Replace-Text $cur "In this narrative, we delve into the captivating realm of quantum biology, exploring its nascent applications in medicine, energy, and information processing" "At the heart of civics lies the exploration of power dynamics and the intricate interplay between various stakeholders"

$cur = Replace-Text $cur " We traverse the intricacies of quantum entanglement in photosynthesis, uncovering how plants harness sunlight with astounding efficiency" " We will examine the concept of sovereignty, analyzing who holds the ultimate authority within a political system and how power is distributed among different entities"

$cur = Delete-Text $cur " We ponder the potential of quantum computing to accelerate drug discovery and protein folding simulations, ushering in an era of personalized medicine and disease mitigation."

$cur = Replace-Text $cur " As we push the boundaries of quantum technology, we unearth new avenues for combating antibiotic resistance, developing efficient energy conversion systems, and navigating the uncharted frontiers of quantum information processing" " Moreover, we will investigate the diverse mechanisms through which citizens can exercise their influence, ranging from voting and lobbying to peaceful protests and community activism"

# =============================================================================
# "Summary" heading: drop the cached <w:lastRenderedPageBreak/> without
# disturbing its run formatting (Aptos / black / 28 half-points).
# =============================================================================
$r = Step-Find $cur "Summary"
$savedFont = $r.Font.Name
$savedSize = $r.Font.Size
$savedColor = $r.Font.Color
$r.Delete()
$r.InsertAfter("Summary")
$r.Font.Name = $savedFont
$r.Font.Size = $savedSize
$r.Font.Color = $savedColor
$cur = $r
$cur.Collapse(0)

# =============================================================================
# Summary paragraph
# =============================================================================
$cur = Replace-Text $cur "Quantum technology has emerged as a vibrant frontier, offering transformative potential in biology" "This essay has provided a comprehensive introduction to the realm of civics, exploring the foundational principles, structures, and processes that underpin our governance systems"

$cur = Replace-Text $cur " Delving into nature's quantum secrets promises a deeper understanding of life's intricate mechanisms, with wide-ranging implications for medicine, materials, energy, and information processing" " We have journeyed through the historical evolution of governance, delved into the intricate interplay of power dynamics, and examined the avenues for citizen engagement"

$cur = Replace-Text $cur " The quest to harness quantum phenomena in biology faces challenges, but the potential rewards are immense" " Through "

$cur = Delete-Text $cur "."

$cur = Replace-Text $cur " From deciphering cellular processes to unraveling DNA replication, quantum biology has the potential to revolutionize our comprehension of life itself" "this exploration, we have gained a deeper understanding of the significance of informed and active citizenship, recognizing the crucial role it plays in shaping our political landscapes"

$cur = Delete-Text $cur "."

$cur = Delete-Text $cur " As we continue to probe the enigmatic quantum world, we stand on the cusp of a new era, where the boundaries of science and technology converge to illuminate the mysteries of life"

# Final trailing period stays untouched.

# =============================================================================
# New empty paragraph at the very end of the body (before sectPr).
# =============================================================================
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
